$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2 with the new zip code value
$ws.Range("A2").Value = 201301

# Remove the rows that are no longer needed (A3:A4), shrinking the used range
$ws.Range("A3:A4").EntireRow.Delete()

# Move the active selection to A3 (now just past the used data)
$ws.Range("A3").Select()
